$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) format, used to force price-like
# strings to remain text instead of being auto-parsed as numbers by Excel.
# C5 (a URL cell) is never touched by this script, so its style stays the
# original default style throughout.
$refStyle = $ws.Range("C5").Style

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $refStyle
}

$ws.Range('D2').Value = '62.867.12'
$ws.Range('D3').Value = '3.440.94'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue 'D5' '576.84'
$ws.Range('E5').Value = '  -0.99%  '
Set-TextValue 'D6' '146.56'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '3.441.32'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('E12').Value = '  +2.71%  '
$ws.Range('D13').Value = '4.025.19'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('E14').Value = '  +2.50%  '
Set-TextValue 'D15' '28.84'
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').Value = '3.421.79'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '62.890.10'
$ws.Range('E18').Value = '  +0.09%  '
Set-TextValue 'D19' '6.34'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('E21').Value = '  -1.48%  '
Set-TextValue 'D22' '384.89'
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('E23').Value = '  -0.49%  '
Set-TextValue 'D24' '74.37'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').Value = '3.578.11'
$ws.Range('E27').Value = '  -3.87%  '
$ws.Range('E28').Value = '  -6.15%  '
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('E30').Value = '  +0.07%  '
Set-TextValue 'D31' '8.07'
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('E32').Value = '  -2.04%  '
$ws.Range('E33').Value = '  -0.04%  '
Set-TextValue 'D34' '23.23'
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('E35').Value = '  -9.83%  '
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('E37').Value = '  -0.26%  '
Set-TextValue 'D38' '31.76'
$ws.Range('E38').Value = '  +3.38%  '
Set-TextValue 'D39' '1.58'
$ws.Range('E39').Value = '  -0.54%  '
Set-TextValue 'D40' '168.38'
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('D41').Value = '3.476.86'
Set-TextValue 'D42' '0.0767'
$ws.Range('E42').Value = '  +0.14%  '
Set-TextValue 'D43' '0.788'
$ws.Range('E43').Value = '  -0.33%  '
Set-TextValue 'D44' '42.35'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('E47').Value = '  -3.24%  '
$ws.Range('D48').Value = '2.565.43'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('E49').Value = '  +2.89%  '
Set-TextValue 'D50' '6.82'
$ws.Range('E50').Value = '  +1.23%  '
Set-TextValue 'D51' '22.56'
$ws.Range('E51').Value = '  -4.07%  '
